$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 935.290659859767
$ws.Range("C2").Value = 949.169319930612
$ws.Range("D2").Value = 962.7163197101
$ws.Range("B3").Value = 568.741523200579
$ws.Range("C3").Value = 577.24063493066
$ws.Range("D3").Value = 585.503102937956
$ws.Range("B4").Value = 405.449306488627
$ws.Range("C4").Value = 411.817135020407
$ws.Range("D4").Value = 418.124684501686
$ws.Range("B5").Value = 596.878222110257
$ws.Range("C5").Value = 606.125156226542
$ws.Range("D5").Value = 615.300270160788
$ws.Range("B6").Value = 1550.36215238271
$ws.Range("C6").Value = 1574.42134302359
$ws.Range("D6").Value = 1598.57951744317
$ws.Range("B7").Value = 1015.33696272899
$ws.Range("C7").Value = 1030.93628594516
$ws.Range("D7").Value = 1046.45769364481
$ws.Range("B8").Value = 1038.24408157393
$ws.Range("C8").Value = 1054.05502725735
$ws.Range("D8").Value = 1069.78613018249
$ws.Range("B9").Value = 570.65040444965
$ws.Range("C9").Value = 579.295446297309
$ws.Range("D9").Value = 587.817459887338
$ws.Range("B10").Value = 506.73124574577
$ws.Range("C10").Value = 513.831042510357
$ws.Range("D10").Value = 520.797065883371
$ws.Range("B11").Value = 1339.97480409885
$ws.Range("C11").Value = 1358.65567204488
$ws.Range("D11").Value = 1376.88614226082
$ws.Range("B12").Value = 1733.40878562539
$ws.Range("C12").Value = 1757.22872000623
$ws.Range("D12").Value = 1780.37500308953
$ws.Range("B13").Value = 1711.50828778401
$ws.Range("C13").Value = 1735.56134366754
$ws.Range("D13").Value = 1759.5505181623
$ws.Range("B14").Value = 725.783565027881
$ws.Range("C14").Value = 735.966090869871
$ws.Range("D14").Value = 745.981402973445
$ws.Range("B15").Value = 1582.97567315157
$ws.Range("C15").Value = 1605.85287299653
$ws.Range("D15").Value = 1628.85458982247
$ws.Range("B16").Value = 363.104338666042
$ws.Range("D16").Value = 410.37950110092
$ws.Range("B17").Value = 361.192176919612
$ws.Range("D17").Value = 408.017517394318
$ws.Range("B18").Value = 344.078450239495
$ws.Range("C18").Value = 348.703558933792
$ws.Range("D18").Value = 353.296371475418
$ws.Range("B19").Value = 338.292076804077
$ws.Range("C19").Value = 342.834142773817
$ws.Range("D19").Value = 347.349050052238
$ws.Range("B20").Value = 1376.56165587883
$ws.Range("C20").Value = 1395.1302780526
$ws.Range("D20").Value = 1412.94255690797
$ws.Range("B21").Value = 1248.63776891335
$ws.Range("D21").Value = 1403.15897493905
$ws.Range("B22").Value = 1233.28110469582
$ws.Range("D22").Value = 1386.97917326292
$ws.Range("B23").Value = 1434.99123686254
$ws.Range("D23").Value = 1613.45865931237
$ws.Range("B24").Value = 1429.82896658252
$ws.Range("D24").Value = 1610.76678453751
$ws.Range("B25").Value = 1822.41332470373
$ws.Range("D25").Value = 2049.37076748755
$ws.Range("B26").Value = 1376.17814789472
$ws.Range("D26").Value = 1546.77875310165
$ws.Range("B27").Value = 1186.84177505301
$ws.Range("D27").Value = 1334.36353390158
$ws.Range("B28").Value = 2456.61564790846
$ws.Range("D28").Value = 2764.24264881697
$ws.Range("B29").Value = 2441.74301581139
$ws.Range("D29").Value = 2747.69901887352
$ws.Range("B30").Value = 2114.03843865272
$ws.Range("D30").Value = 2377.00712085549
$ws.Range("B31").Value = 2052.81004074937
$ws.Range("D31").Value = 2309.53984380317
$ws.Range("B32").Value = 482.911910649433
$ws.Range("D32").Value = 544.071426909266
$ws.Range("B33").Value = 518.207769602502
$ws.Range("D33").Value = 584.43443412136
$ws.Range("B34").Value = 622.976314845446
$ws.Range("C34").Value = 631.526849583293
$ws.Range("D34").Value = 639.863940396487
$ws.Range("B35").Value = 1296.7410188606
$ws.Range("C35").Value = 1314.28484770269
$ws.Range("D35").Value = 1331.04306246265
$ws.Range("B36").Value = 1048.92930777992
$ws.Range("C36").Value = 1063.1121073411
$ws.Range("D36").Value = 1076.78196780946
$ws.Range("B37").Value = 662.485296933878
$ws.Range("C37").Value = 671.44831876705
$ws.Range("D37").Value = 680.238923676199
$ws.Range("B38").Value = 1238.84515417432
$ws.Range("C38").Value = 1255.69241748671
$ws.Range("D38").Value = 1271.92305723476
$ws.Range("B39").Value = 1274.97115761175
$ws.Range("C39").Value = 1292.39824074298
$ws.Range("D39").Value = 1309.19394878223
$ws.Range("B40").Value = 1060.590484961
$ws.Range("C40").Value = 1075.05063274154
$ws.Range("D40").Value = 1088.94904150001
$ws.Range("B41").Value = 1153.89012402349
$ws.Range("C41").Value = 1169.74135892255
$ws.Range("D41").Value = 1184.9359618325
$ws.Range("B42").Value = 1155.88669322891
$ws.Range("C42").Value = 1171.68079542622
$ws.Range("D42").Value = 1187.01288409883
$ws.Range("B43").Value = 475.938908636944
$ws.Range("C43").Value = 482.592326589112
$ws.Range("D43").Value = 489.115085178555
$ws.Range("B44").Value = 985.323396438844
$ws.Range("C44").Value = 998.816007626036
$ws.Range("D44").Value = 1011.80621379709
$ws.Range("B45").Value = 1535.70254052128
$ws.Range("C45").Value = 1557.56182306632
$ws.Range("D45").Value = 1579.39986834752
$ws.Range("B46").Value = 1528.77883144404
$ws.Range("C46").Value = 1550.49528063082
$ws.Range("D46").Value = 1572.25264980169
$ws.Range("B47").Value = 639.99697247528
$ws.Range("C47").Value = 648.947715757392
$ws.Range("D47").Value = 657.73179010961
$ws.Range("B48").Value = 536.948656657408
$ws.Range("C48").Value = 544.176543723348
$ws.Range("D48").Value = 551.284137517882
